$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2009
$ws.Range("D3").Value = 2006
$ws.Range("D4").Value = 2000
$ws.Range("D5").Value = 1998
$ws.Range("D6").Value = 1995
$ws.Range("D7").Value = 1993
$ws.Range("D8").Value = 1992
$ws.Range("D9").Value = 1991
$ws.Range("D10").Value = 1990

$ws.Range("D1").Select()
